$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1954.8889
$ws.Range("I62").Value = 1822.3077
$ws.Range("J62").Value = 2299.6
$ws.Range("K62").Value = 1822.3077
$ws.Range("L62").Value = 2299.6
$ws.Range("M62").Value = -1198.3077
$ws.Range("N62").Value = -3547.6
$ws.Range("H64").Value = 3656.9644
$ws.Range("I64").Value = 3845
$ws.Range("J64").Value = 3186.875
$ws.Range("K64").Value = 3845
$ws.Range("L64").Value = 3186.875
$ws.Range("M64").Value = -3597
$ws.Range("N64").Value = -3682.875
$ws.Range("H65").Value = 1954.8889
$ws.Range("I65").Value = 1822.3077
$ws.Range("J65").Value = 2299.6
$ws.Range("K65").Value = 9111.538500000001
$ws.Range("L65").Value = 11498
$ws.Range("M65").Value = -5991.538500000001
$ws.Range("N65").Value = -17738
$ws.Range("H67").Value = 3656.9644
$ws.Range("I67").Value = 3845
$ws.Range("J67").Value = 3186.875
$ws.Range("K67").Value = 3845
$ws.Range("L67").Value = 3186.875
$ws.Range("M67").Value = -2987
$ws.Range("N67").Value = -4902.875
$ws.Range("H69").Value = 3294106
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3294106
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9882318
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -9884066
$ws.Range("H72").Value = 3294106
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3294106
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 29646954
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -29655690
$ws.Range("H74").Value = 3843.3333
$ws.Range("I74").Value = 3843.3333
$ws.Range("K74").Value = 3843.3333
$ws.Range("M74").Value = -2907.3333
$ws.Range("H76").Value = 6344.4736
$ws.Range("I76").Value = 6999.4614
$ws.Range("J76").Value = 4925.3335
$ws.Range("K76").Value = 6999.4614
$ws.Range("L76").Value = 4925.3335
$ws.Range("M76").Value = -6684.4614
$ws.Range("N76").Value = -5555.3335
$ws.Range("H77").Value = 3843.3333
$ws.Range("I77").Value = 3843.3333
$ws.Range("K77").Value = 19216.6665
$ws.Range("M77").Value = -14536.6665
$ws.Range("H79").Value = 6344.4736
$ws.Range("I79").Value = 6999.4614
$ws.Range("J79").Value = 4925.3335
$ws.Range("K79").Value = 6999.4614
$ws.Range("L79").Value = 4925.3335
$ws.Range("M79").Value = -5907.4614
$ws.Range("N79").Value = -7109.3335
$ws.Range("H80").Value = 230.57143
$ws.Range("J80").Value = 170
$ws.Range("L80").Value = 510
$ws.Range("N80").Value = -2506
$ws.Range("H83").Value = 230.57143
$ws.Range("J83").Value = 170
$ws.Range("L83").Value = 1530
$ws.Range("N83").Value = -11514
$ws.Range("H100").Value = 15153325
$ws.Range("I100").Value = 27779496
$ws.Range("J100").Value = 1920
$ws.Range("K100").Value = 27779496
$ws.Range("L100").Value = 1920
$ws.Range("M100").Value = -27778955
$ws.Range("N100").Value = -3002
$ws.Range("H129").Value = 1010.38574
$ws.Range("I129").Value = 386
$ws.Range("J129").Value = 1102.5082
$ws.Range("K129").Value = 1158
$ws.Range("L129").Value = 3307.5246
$ws.Range("M129").Value = 3842
$ws.Range("N129").Value = -13307.5246
$ws.Range("H132").Value = 904.79364
$ws.Range("I132").Value = 812.06665
$ws.Range("K132").Value = 2436.19995
$ws.Range("M132").Value = 93.80004999999983
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6712.567
$ws.Range("I32").Value = 4733.3374
$ws.Range("K32").Value = 4733.3374
$ws.Range("M32").Value = -4446.3374
$ws.Range("H61").Value = 235398.67
$ws.Range("I61").Value = 2100.1738
$ws.Range("J61").Value = 503691.94
$ws.Range("K61").Value = 2100.1738
$ws.Range("L61").Value = 503691.94
$ws.Range("M61").Value = -1888.1738
$ws.Range("N61").Value = -504115.94
$ws.Range("H74").Value = 1376.2188
$ws.Range("I74").Value = 1175.5217
$ws.Range("J74").Value = 1889.1111
$ws.Range("K74").Value = 1175.5217
$ws.Range("L74").Value = 1889.1111
$ws.Range("M74").Value = -301.5217
$ws.Range("N74").Value = -3637.1111
$ws.Range("H77").Value = 1376.2188
$ws.Range("I77").Value = 1175.5217
$ws.Range("J77").Value = 1889.1111
$ws.Range("K77").Value = 5877.6085
$ws.Range("L77").Value = 9445.5555
$ws.Range("M77").Value = -1509.6085
$ws.Range("N77").Value = -18181.5555
$ws.Range("H136").Value = 235398.67
$ws.Range("I136").Value = 2100.1738
$ws.Range("J136").Value = 503691.94
$ws.Range("K136").Value = 6300.5214
$ws.Range("L136").Value = 1511075.82
$ws.Range("M136").Value = -3750.5214
$ws.Range("N136").Value = -1516175.82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2030.4
$ws.Range("I86").Value = 1906.2222
$ws.Range("J86").Value = 2216.6667
$ws.Range("K86").Value = 1906.2222
$ws.Range("L86").Value = 2216.6667
$ws.Range("M86").Value = -783.2221999999999
$ws.Range("N86").Value = -4462.6667
$ws.Range("H89").Value = 2030.4
$ws.Range("I89").Value = 1906.2222
$ws.Range("J89").Value = 2216.6667
$ws.Range("K89").Value = 9531.110999999999
$ws.Range("L89").Value = 11083.3335
$ws.Range("M89").Value = -3915.110999999999
$ws.Range("N89").Value = -22315.3335
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6184.875
$ws.Range("I62").Value = 11240
$ws.Range("J62").Value = 4499.8335
$ws.Range("K62").Value = 11240
$ws.Range("L62").Value = 4499.8335
$ws.Range("M62").Value = -10616
$ws.Range("N62").Value = -5747.8335
$ws.Range("H65").Value = 6184.875
$ws.Range("I65").Value = 11240
$ws.Range("J65").Value = 4499.8335
$ws.Range("K65").Value = 56200
$ws.Range("L65").Value = 22499.1675
$ws.Range("M65").Value = -53080
$ws.Range("N65").Value = -28739.1675
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6349.3213
$ws.Range("I80").Value = 9053.666999999999
$ws.Range("J80").Value = 3228.923
$ws.Range("K80").Value = 9053.666999999999
$ws.Range("L80").Value = 3228.923
$ws.Range("M80").Value = -8055.666999999999
$ws.Range("N80").Value = -5224.923
$ws.Range("H83").Value = 6349.3213
$ws.Range("I83").Value = 9053.666999999999
$ws.Range("J83").Value = 3228.923
$ws.Range("K83").Value = 45268.335
$ws.Range("L83").Value = 16144.615
$ws.Range("M83").Value = -40276.335
$ws.Range("N83").Value = -26128.615
$ws.Range("H93").Value = 9251
$ws.Range("J93").Value = 9251
$ws.Range("L93").Value = 9251
$ws.Range("N93").Value = -12995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1380272.4
$ws.Range("I82").Value = 2000636
$ws.Range("J82").Value = 346333
$ws.Range("K82").Value = 2000636
$ws.Range("L82").Value = 346333
$ws.Range("M82").Value = -2000275
$ws.Range("N82").Value = -347055
$ws.Range("H85").Value = 1380272.4
$ws.Range("I85").Value = 2000636
$ws.Range("J85").Value = 346333
$ws.Range("K85").Value = 2000636
$ws.Range("L85").Value = 346333
$ws.Range("M85").Value = -1999388
$ws.Range("N85").Value = -348829
$ws.Range("H132").Value = 12823969
$ws.Range("I132").Value = 20836020
$ws.Range("J132").Value = 4687.8
$ws.Range("K132").Value = 62508060
$ws.Range("L132").Value = 14063.4
$ws.Range("M132").Value = -62505530
$ws.Range("N132").Value = -19123.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2615.1052
$ws.Range("I136").Value = 3154.6667
$ws.Range("J136").Value = 2129.5
$ws.Range("K136").Value = 9464.000100000001
$ws.Range("L136").Value = 6388.5
$ws.Range("M136").Value = -6914.000100000001
$ws.Range("N136").Value = -11488.5
